# Update row 9 (Ano = 2025) with refreshed "Dados BIBI" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3571715.55
$ws.Range("C9").Value = 561545.24
$ws.Range("D9").Value = 4133260.79
$ws.Range("E9").Value = 13.58601038092252
$ws.Range("F9").Value = 86.41398961907748
$ws.Range("G9").Value = -45.72945254765099
$ws.Range("H9").Value = -35.49975413372506
$ws.Range("I9").Value = 35958
$ws.Range("J9").Value = 1536
$ws.Range("K9").Value = 37494
$ws.Range("L9").Value = 25887
$ws.Range("M9").Value = 159.6654996716499
$ws.Range("N9").Value = 9.006888649528744
